$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-classified MAGs (rows 2-9), columns A-L: name, 8 genus
# probabilities, max, prediction, rejection-f.
# Built as a true 2D array since Range.Value2 requires a rectangular
# array (a jagged array of row-arrays is not accepted).
$rowNames = @(
    "even_MAG-GUT25429.fa"
    "even_MAG-GUT30539.fa"
    "even_MAG-GUT32544.fa"
    "even_MAG-GUT44190.fa"
    "even_MAG-GUT44598.fa"
    "even_MAG-GUT56457.fa"
    "even_MAG-GUT58695.fa"
    "even_MAG-GUT66161.fa"
)

$rowTails = @(
    ,([double]"0.0001161187098641203", [double]"0.3375537752363922", [double]"8.876392087183043e-06", [double]"4.632730769288377e-05", [double]"0.0004247387306391607", [double]"0.0004312413320542567", [double]"0.002757412263511772", [double]"0.6586615100277585", [double]"0.6586615100277585", "g__Terrisporobacter", "g__Terrisporobacter")
    ,([double]"0.004988653267735906", [double]"0.1442081357704421", [double]"0.0003381728358450314", [double]"0.0003917443835891789", [double]"0.004579497887763636", [double]"0.006587718642567897", [double]"0.01162078558482822", [double]"0.8272852916272281", [double]"0.8272852916272281", "g__Terrisporobacter", "g__Terrisporobacter")
    ,([double]"4.04302555443032e-05", [double]"0.5835253975811087", [double]"4.333772000862312e-07", [double]"5.588968055731193e-05", [double]"0.0002098667123155597", [double]"0.0001094422141574709", [double]"0.007194975951057518", [double]"0.4088635642280593", [double]"0.5835253975811087", "g__GCA-900066495", "g__GCA-900066495")
    ,([double]"0.001043055964462028", [double]"0.1065710652265113", [double]"0.000698215239215151", [double]"0.0002652024919012892", [double]"0.002998658823059914", [double]"0.01107359323701635", [double]"0.00527409989671241", [double]"0.8720761091211215", [double]"0.8720761091211215", "g__Terrisporobacter", "g__Terrisporobacter")
    ,([double]"0.02093453778093218", [double]"0.2365663271633392", [double]"0.05833908912843729", [double]"0.06079360040260021", [double]"0.01914924809416653", [double]"0.0002884440280362653", [double]"0.531598451611504", [double]"0.07233030179098429", [double]"0.531598451611504", "g__Romboutsia", "g__Romboutsia")
    ,([double]"0.0008024529210328089", [double]"0.2670278290007201", [double]"9.554374094889125e-05", [double]"2.568700667197349e-05", [double]"0.002353378892070749", [double]"0.004701166315230735", [double]"0.003759029129442887", [double]"0.7212349129938819", [double]"0.7212349129938819", "g__Terrisporobacter", "g__Terrisporobacter")
    ,([double]"0.007672474862095388", [double]"0.07980807694049", [double]"0.001185724538680015", [double]"8.014235894879851e-05", [double]"0.001627950153799406", [double]"0.01290616761492669", [double]"0.0009432756091127345", [double]"0.895776187921947", [double]"0.895776187921947", "g__Terrisporobacter", "g__Terrisporobacter")
    ,([double]"0.0007669148399584079", [double]"0.5213827587689047", [double]"8.828638910138277e-06", [double]"0.00013671675215724", [double]"0.0008917292209083694", [double]"0.0004742692527444682", [double]"0.008987803119281761", [double]"0.467350979407135", [double]"0.5213827587689047", "g__GCA-900066495", "g__GCA-900066495")
)

$data = New-Object "object[,]" 8,12
for ($i = 0; $i -lt 8; $i++) {
    $data[$i,0] = $rowNames[$i]
    $tail = $rowTails[$i]
    for ($j = 0; $j -lt 11; $j++) {
        $data[$i, $j + 1] = $tail[$j]
    }
}

$rng = $ws.Range("A2:L9")
$rng.Value2 = $data

# Column A keeps the bordered/centered style used by the header row and original A2;
# copy-format (rather than `.Style =`) so the existing cellXfs entry is reused instead
# of a new (merely equivalent) style being created.
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
